# Add max root depth msmts
# - rename "stage_Vx" header to "stage"
# - column G values (stage) become text labels "V12"/"V14"/"V15" instead of
#   the bare numeric VN stage number (column H keeps the numeric value)
# - selection/scroll position updated to G7

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header rename: G6 "stage_Vx" -> "stage"
$ws.Range("G6").Value = "stage"

# Map of row -> new text label for column G
$gValues = @{
    7  = "V14"
    8  = "V15"
    9  = "V14"
    10 = "V14"
    11 = "V14"
    12 = "V14"
    13 = "V14"
    14 = "V15"
    15 = "V14"
    16 = "V14"
    17 = "V15"
    18 = "V14"
    19 = "V15"
    20 = "V14"
    21 = "V14"
    22 = "V12"
    23 = "V12"
    24 = "V12"
    25 = "V14"
    26 = "V14"
    27 = "V15"
    28 = "V15"
    29 = "V14"
    30 = "V14"
}

foreach ($row in $gValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $gValues[$row]
}

# A handful of rows (20, 21, 25-30) used the larger/bold-ish "summary row"
# font+alignment (copied down from the totals style) for column G; once
# retyped as text they match the plain centered style used elsewhere in
# the column (e.g. G10: default Calibri 11, theme text color, centered
# both ways), so re-format just those rows to match.
$restyleRows = @(20, 21, 25, 26, 27, 28, 29, 30)
foreach ($row in $restyleRows) {
    $cell = $ws.Cells.Item($row, 7)
    $cell.Font.Name = "Calibri"
    $cell.Font.Size = 11
    $cell.Font.ThemeColor = 1
    $cell.Font.ThemeFont = 1
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4108
}

# Update selection / scroll position
$ws.Activate() | Out-Null
$ws.Range("G7").Select() | Out-Null
